# Apply the "material" column edit described in the commit:
#   "added primitive version of reading passages and changed question
#    model to have a 'material' field"
#
# This adds a new column H ("has_material") to the questions sheet with a
# yes/no flag per question row, narrows/re-sizes columns D:G (which used
# to be auto "best fit" sized for the old content), and updates the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the new "has_material" column, row by row (row 1 is the header).
$hasMaterial = @{
    1  = "has_material"
    2  = "yes"
    3  = "no"
    4  = "no"
    5  = "no"
    6  = "no"
    7  = "yes"
    8  = "no"
    9  = "no"
    10 = "no"
    11 = "no"
    12 = "no"
    13 = "yes"
    14 = "yes"
    15 = "no"
    16 = "no"
    17 = "no"
    18 = "no"
    19 = "no"
    20 = "no"
    21 = "no"
}

foreach ($r in 1..21) {
    $ws.Cells.Item($r, 8).Value = $hasMaterial[$r]
}

# Row 21's answer cells use the wrap-text style; match that on the new cell too.
$ws.Cells.Item(21, 8).WrapText = $true

# The old columns D:G were auto "best fit"; now that there's a narrower,
# fixed-purpose "material" column they get fixed, smaller widths instead.
$ws.Columns.Item(4).ColumnWidth = 21.85546875
$ws.Columns.Item(5).ColumnWidth = 20.140625
$ws.Columns.Item(6).ColumnWidth = 19.7109375
$ws.Columns.Item(7).ColumnWidth = 24.5703125

# Row 21 grew taller once the columns reflowed.
$ws.Rows.Item(21).RowHeight = 105

# Update the active selection to reflect where the user was last working.
$ws.Range("C13").Select() | Out-Null
